# Fruta / hortaliza, semanal
#
# Insert a new weekly price record for "Piña" (Feria Lagunitas de Puerto
# Montt) just above the existing row 368, pushing all subsequent rows
# down by one (old row 396 becomes row 397). The sheet's used range grows
# from A1:T396 to A1:T397 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 368 (and everything below it) down by one row.
$ws.Rows.Item(368).Insert()

# Populate the newly-inserted row with the new record's data.
$ws.Cells.Item(368, 1).Value  = 4
$ws.Cells.Item(368, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(368, 3).Value  = "Los Lagos"
$ws.Cells.Item(368, 4).Value  = 45013
$ws.Cells.Item(368, 5).Value  = 10
$ws.Cells.Item(368, 6).Value  = "Fruta"
$ws.Cells.Item(368, 7).Value  = 100108
$ws.Cells.Item(368, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(368, 9).Value  = 100108005
$ws.Cells.Item(368, 10).Value = "Piña"
$ws.Cells.Item(368, 11).Value = "Caramelo"
$ws.Cells.Item(368, 12).Value = "Primera"
$ws.Cells.Item(368, 13).Value = 200
$ws.Cells.Item(368, 14).Value = 21000
$ws.Cells.Item(368, 15).Value = 22000
$ws.Cells.Item(368, 16).Value = 21500
$ws.Cells.Item(368, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(368, 18).Value = "Ecuador"
$ws.Cells.Item(368, 19).Value = 1792
$ws.Cells.Item(368, 20).Value = 12
